$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new record row (row 5) with the same layout/formatting as the
# existing rows: Nombre, Apellido, Edad, Fecha, Hora
$ws.Range("A5").Value = "Kevin"
$ws.Range("B5").Value = "Nu" + [char]0x00F1 + "ez"
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 44732
$ws.Range("E5").Value = 0.86458333333333337

# Copy formatting from the row above so the new row matches existing style
$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A5:E5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the values (paste of formats only should not disturb them, but
# make sure number formats didn't overwrite our intended values)
$ws.Range("A5").Value = "Kevin"
$ws.Range("B5").Value = "Nu" + [char]0x00F1 + "ez"
$ws.Range("C5").Value = 21
$ws.Range("D5").Value = 44732
$ws.Range("E5").Value = 0.86458333333333337

# Update the active selection to match the saved view state
$ws.Range("F14").Select() | Out-Null

$wb.Save()
